$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row 3 (shifts existing rows 3-11 down to 4-12)
$ws.Rows.Item(3).Insert()

# Populate the new row with the dataModelType / solution pair
$ws.Cells.Item(3, 1).Value = "dataModelType"
$ws.Cells.Item(3, 2).Value = "solution"

# The hyperlink that used to live on B5 now lives on B6 (content shifted down
# by the row insert); re-create it pointing at the new location so the
# worksheet XML records the correct <hyperlink ref="B6" .../> entry.
$hyperlinkCell = $ws.Cells.Item(6, 2)
$hyperlinkTarget = "http://purl.org/cognite/power_analytic"
$originalStyle = $hyperlinkCell.Style
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($hyperlinkCell, $hyperlinkTarget)
$hyperlinkCell.Style = $originalStyle

# Update the selection to match the author's final cursor position
$ws.Range("B26").Select()
